# edit.ps1 - apply the "updated quarto and minor edits" change set
#
# Summary of changes (see commit diff):
#   1. Collapse the "N. " run + heading-text run (in each of the 4 table
#      header cells) into a single run with the combined text.
#   2. In the "Pro tip" paragraph, change " and collaborate on " -> " and "
#      and "others" -> "work on other tasks together", while keeping those
#      two runs *separate* (not merged into one run).
#   3. Best-effort: mark the built-in "Default Paragraph Font" character
#      style as semi-hidden (w:semiHidden) to match the stylesheet tweak
#      in the diff. (Left wrapped in try/catch: this runtime's Style
#      object does not expose a writable Hidden/SemiHidden property, so
#      this is attempted but allowed to silently no-op if unsupported.)

$d = $word.ActiveDocument

# --- 1) Merge "N. " + heading text runs into a single run per cell -------
# Find.Execute can match text that spans multiple runs; when it replaces
# the matched span it collapses the runs it touched into one run carrying
# the full replacement text (and drops the now-superfluous rsid attrs),
# which is exactly the shape the target XML has. Since we want the same
# text back (just merged), replace each heading with itself.

$d.Content.Find.Execute(
    "1. Data and Variables Description", $true, $false, $false, $false,
    $false, $true, 1, $false,
    "1. Data and Variables Description", 2) | Out-Null

$d.Content.Find.Execute(
    "2. Summary Statistics", $true, $false, $false, $false,
    $false, $true, 1, $false,
    "2. Summary Statistics", 2) | Out-Null

$d.Content.Find.Execute(
    "3. Relationship b/w Main Variables", $true, $false, $false, $false,
    $false, $true, 1, $false,
    "3. Relationship b/w Main Variables", 2) | Out-Null

$d.Content.Find.Execute(
    "4. Exploring Other Variables", $true, $false, $false, $false,
    $false, $true, 1, $false,
    "4. Exploring Other Variables", 2) | Out-Null

# --- 2) "... and collaborate on others." -> "... and work on other tasks together." ---
# Editing text in-place (via Find or Range.Text) causes this runtime to
# coalesce every contiguous run that ends up sharing identical run
# formatting (<w:rPr>) with the edited run, on both sides. The target XML
# keeps " and " and "work on other tasks together" as two distinct runs,
# so each edit is done on a Range whose formatting is briefly nudged
# (Bold toggled on, then back off) so it does not get folded into its
# (same-format) neighbours "tasks" / "others"/".". Toggling the format
# back off afterwards does not retroactively re-trigger a merge.
#
# Find.Execute is run directly against a Range object (not a freshly
# fetched $d.Content each time); on a match that Range's Start/End/Text
# collapse onto the found text, which gives us an exact, reliable anchor
# for the subsequent Range-based edit (wdFindStop / no replacement here).

$needle1 = " and collaborate on "
$r1 = $d.Content
$found1 = $r1.Find.Execute(
    $needle1, $true, $false, $false, $false, $false, $true, 1, $false,
    "", 0)

$r1.Bold = 1
$r1.Text = " and "
$r1.Bold = 0

$r2 = $d.Range($r1.End, $r1.End + 6)
$r2.Bold = 1
$r2.Text = "work on other tasks together"
$r2.Bold = 0

# --- 3) Best-effort: styles.xml "Default Paragraph Font" semiHidden -----
try {
    $dpf = $d.Styles("Default Paragraph Font")
    $dpf.Hidden = $true
} catch {
    # Not supported by this host's Style object model surface; skip.
}
